# Add two new fields ("celebrity" and "brands") to the header row of the
# YouTube data import template, as per commit message:
# "Updated template as per celebrity and brand field added"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells appended right after the existing last column (Y -> Z, AA)
$ws.Cells.Item(1, 26).Value = "celebrity"
$ws.Cells.Item(1, 27).Value = "brands"
